$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update B3 (Tareas/total value) from 185 to 194 for hito 3
$ws.Range("B3").Value = 194

# Recalculate formulas so dependent cells (B7, G7, I9) reflect the new value
$excel.Calculate()
